$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS) - update VENTA and POR CUMPLIR
$ws.Range("D2").Value = 925.91
$ws.Range("E2").Value = -925.91

# Row 4 (TOTAL) - update VENTA, POR CUMPLIR and CUMPLIMIENTO
$ws.Range("D4").Value = 1257.66
$ws.Range("E4").Value = 12465.68
$ws.Range("F4").Value = 0.09164387095269809
